# Applies the "Collection of FlexTest 20mm" commit: 17 data points for the
# flexor muscle (20mm BPA) plus a couple of small corrections on the
# ExtTest40mm sheet.

$wb = $excel.ActiveWorkbook
$wsExt = $wb.Worksheets.Item("ExtTest40mm")
$wsFlx = $wb.Worksheets.Item("FlxTest20mm")

# ---------------------------------------------------------------------
# FlxTest20mm ("flexor 20mm") sheet — new columns D:S (tests 2-17), a new
# "tendon length" block in row 4, and two footnote/annotation rows.
# Shared strings must be introduced in this exact order so they line up
# with the target sharedStrings.xml (index 18..24).
# ---------------------------------------------------------------------

# -- C2 / C3 (overall numbers) ---------------------------------------
$wsFlx.Range("C2").Value = 423
# C3 formula is unchanged (=C2-C2*0.17); recalculated automatically.

# -- new M3 computed cell ---------------------------------------------
$wsFlx.Range("M3").Formula = "=1-369.5/423"

# -- new row 4: "tendon length" + footnote markers --------------------
$wsFlx.Range("B4").Value = "tendon length"   # shared string 18
$wsFlx.Range("C4").Value = 17
$wsFlx.Range("I4").Value = "*"               # shared string 19
$wsFlx.Range("M4").Value = "**"              # shared string 20

# -- row 5: extend test-number header to tests 16 & 17 -----------------
$wsFlx.Range("R5").Value = 16
$wsFlx.Range("S5").Value = 17

# -- row 6: Load (N) ----------------------------------------------------
$wsFlx.Range("C6").Value = 29.225
$wsFlx.Range("D6").Value = 24.009
$wsFlx.Range("E6").Value = 20.302
$wsFlx.Range("F6").Value = 17.246
$wsFlx.Range("G6").Value = 14.472
$wsFlx.Range("H6").Value = 11.127
$wsFlx.Range("I6").Value = 11.121
$wsFlx.Range("J6").Value = 9.284
$wsFlx.Range("K6").Value = 6.8421
$wsFlx.Range("L6").Value = 2.9257
$wsFlx.Range("M6").Value = 19.083
$wsFlx.Range("N6").Value = 14.757
$wsFlx.Range("O6").Value = 11.612
$wsFlx.Range("P6").Value = 5.8364
$wsFlx.Range("Q6").Value = 3.1772
$wsFlx.Range("R6").Value = 7.6226
$wsFlx.Range("S6").Value = 4.8666

# -- row 7: Knee angle ----------------------------------------------------
$wsFlx.Range("D7").Value = 11.5
$wsFlx.Range("E7").Value = 20
$wsFlx.Range("F7").Value = 27
$wsFlx.Range("G7").Value = 35
$wsFlx.Range("H7").Value = 44
$wsFlx.Range("I7").Value = 54.5
$wsFlx.Range("J7").Value = 50
$wsFlx.Range("K7").Value = 63
$wsFlx.Range("L7").Value = 64
$wsFlx.Range("M7").Value = 73.5
$wsFlx.Range("N7").Value = 82
$wsFlx.Range("O7").Value = 88
$wsFlx.Range("P7").Value = 103.5
$wsFlx.Range("Q7").Value = 114
$wsFlx.Range("R7").Value = 115
$wsFlx.Range("S7").Value = 120

# -- row 8: muscle length -------------------------------------------------
$wsFlx.Range("C8").Value = 17.4
$wsFlx.Range("D8").Value = 21.4
$wsFlx.Range("E8").Value = 22.3
$wsFlx.Range("F8").Value = 23.1
$wsFlx.Range("G8").Value = 24.9
$wsFlx.Range("H8").Value = 28.4
$wsFlx.Range("I8").Value = 29.3
$wsFlx.Range("J8").Value = 29.4
$wsFlx.Range("K8").Value = 35.1
$wsFlx.Range("L8").Value = 35.7
$wsFlx.Range("M8").Value = 34.7
$wsFlx.Range("N8").Value = 38.3
$wsFlx.Range("O8").Value = 34.9
$wsFlx.Range("P8").Value = 35.3
$wsFlx.Range("Q8").Value = 42.7
$wsFlx.Range("R8").Value = 46.4
$wsFlx.Range("S8").Value = 48.1

# -- row 9: dl/dtheta -------------------------------------------------
$wsFlx.Range("C9").Value = 24.5
$wsFlx.Range("D9").Value = 22
$wsFlx.Range("E9").Value = 24
$wsFlx.Range("F9").Value = 25.5
$wsFlx.Range("G9").Value = 31.5
$wsFlx.Range("H9").Value = 31
$wsFlx.Range("I9").Value = 31.5
$wsFlx.Range("J9").Value = 30
$wsFlx.Range("K9").Value = 31.5
$wsFlx.Range("L9").Value = 33.5
$wsFlx.Range("M9").Value = 38.5
$wsFlx.Range("N9").Value = 44.5
$wsFlx.Range("O9").Value = 43
$wsFlx.Range("P9").Value = 54.5
$wsFlx.Range("Q9").Value = 55.5
$wsFlx.Range("R9").Value = 52
$wsFlx.Range("S9").Value = 55

# -- row 10: Resting muscle length -------------------------------------
$wsFlx.Range("C10").Value = 401
$wsFlx.Range("D10").Value = 395
$wsFlx.Range("E10").Value = 389
$wsFlx.Range("F10").Value = 387
$wsFlx.Range("G10").Value = 382
$wsFlx.Range("H10").Value = 380
$wsFlx.Range("I10").Value = 378
$wsFlx.Range("J10").Value = 374.5
$wsFlx.Range("K10").Value = 371.5
$wsFlx.Range("L10").Value = 369.5
$wsFlx.Range("M10").Value = 354
$wsFlx.Range("N10").Value = 352.5
$wsFlx.Range("O10").Value = 350.5
$wsFlx.Range("P10").Value = 345
$wsFlx.Range("Q10").Value = 343.5
$wsFlx.Range("R10").Value = 338.5
$wsFlx.Range("S10").Value = 339

# -- row 12: Tibia origin (matlab) -- extend existing 366.43 series -----
$wsFlx.Range("R12").Value = 366.43
$wsFlx.Range("S12").Value = 366.43

# -- row 13: Load cell angle (tibia) -------------------------------------
$wsFlx.Range("C13").Value = 45
$wsFlx.Range("D13").Value = 47
$wsFlx.Range("E13").Value = 51
$wsFlx.Range("F13").Value = 45
$wsFlx.Range("G13").Value = 40
$wsFlx.Range("H13").Value = 37
$wsFlx.Range("I13").Value = 38
$wsFlx.Range("J13").Value = 39
$wsFlx.Range("K13").Value = 36.5
$wsFlx.Range("L13").Value = 35.5
$wsFlx.Range("M13").Value = 30.5
$wsFlx.Range("N13").Value = 20
$wsFlx.Range("O13").Value = 22
$wsFlx.Range("P13").Value = 16
$wsFlx.Range("Q13").Value = 4.5
$wsFlx.Range("R13").Value = 3
$wsFlx.Range("S13").Value = 10

# -- row 15: extend the shared "expected max contract" formula to R:S ---
$wsFlx.Range("R15").Formula = "=R6*COS(RADIANS(R9-2.83))*R12/1000"
$wsFlx.Range("S15").Formula = "=S6*COS(RADIANS(S9-2.83))*S12/1000"

# -- footnote rows 16-18 --------------------------------------------------
$wsFlx.Range("J16").Value = "*Pipe Started to kink at 7"                               # shared string 21
$wsFlx.Range("M17").Value = "**Pressure stepped up from 300kPa to 500kPa on test 11-15" # shared string 22
$wsFlx.Range("O18").Value = "***Pressure stepped up from 500kPa to 620kPa on test 16"   # shared string 23

# -- finish off row 4 with the "***" marker (added last -> shared string 24) --
$wsFlx.Range("R4").Value = "***"

# ---------------------------------------------------------------------
# ExtTest40mm sheet — a handful of corrected readings.
# ---------------------------------------------------------------------
$wsExt.Range("C6").Value = 21.395
$wsExt.Range("C7").Value = 117
$wsExt.Range("C8").Value = 36.1
$wsExt.Range("C9").Value = 32
$wsExt.Range("C10").Value = 533
$wsExt.Range("C13").Value = 60

# ---------------------------------------------------------------------
# View state: ExtTest40mm becomes the active/selected tab with D7
# selected; FlxTest20mm keeps its own last selection (W14) but is no
# longer the active sheet.
# ---------------------------------------------------------------------
$wsFlx.Range("W14").Select() | Out-Null
$wsExt.Activate()
$wsExt.Range("D7").Select() | Out-Null
